# "fix create inv increase"
# The redundant "Chi tiết đơn hàng / Thuế(%)" column is removed from each of
# the three PO sheets (it duplicates the more detailed tax columns further to
# the right: % Thuế nhập khẩu, % Thuế tiêu thụ đặc biệt, % Thuế GTGT, ...).
# Deleting the whole column shifts everything to its right one column to the
# left, and the now-unused shared string is dropped automatically on save.

$wb = $excel.ActiveWorkbook

# Sheet "PO - HH": column T holds "Chi tiết đơn hàng / Thuế(%)"
$wsHH = $wb.Worksheets.Item(1)
$wsHH.Columns("T").Delete()

# Sheet "PO - DV": column Q holds "Chi tiết đơn hàng / Thuế(%)"
$wsDV = $wb.Worksheets.Item(2)

# Fix up the conditional-formatting range that spans the deleted column so it
# keeps referencing the correct (now-narrower) block of cells, the way Excel
# itself would when a column inside a formatted range is removed.
$cfRange = $wsDV.Range("O2:U2")
$cfRule = $cfRange.FormatConditions.Item(1)

$wsDV.Columns("Q").Delete()

$cfRule.ModifyAppliesToRange($wsDV.Range("O2:T2"))
$cfRule.Priority = 5

# Sheet "PO - TS": column Q holds "Chi tiết đơn hàng / Thuế(%)"
$wsTS = $wb.Worksheets.Item(3)
$wsTS.Columns("Q").Delete()
